# Timetable-template.xlsx update:
#  - Replace the "Subject" sheet content with the richer subject table
#    (department/name/hours/preferred slot/capacities/lab type/instructors)
#    that used to live on the separate "NewSubject" sheet.
#  - Drop the now-redundant "NewSubject" sheet.
#  - Extend "Enrollment" with a SUBJECT_ID column (inserted before STUDENT_ID).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Rebuild the "Subject" sheet with the new column layout & data.
# ---------------------------------------------------------------------------
$subject = $wb.Worksheets.Item("Subject")
$subject.Cells.Clear() | Out-Null

$headers = @("DEPARTMENT","NAME","NUM_LAB_HOURS","NUM_LEC_HOURS","PREFERED_WEEKDAY","PREFERED_TIME","LEC_CAPACITY","LAB_CAPACITY","LAB_TYPE","INSTRUCTORS")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $subject.Cells.Item(1, $col + 1).Value = $headers[$col]
}

$subjectRows = @(
    @("MT","GT1",2,3,"0,1,2","7-12,12-18",4,2,"lec","0,1"),
    @("PH","VL1",0,2,"0,1,2","7-12,12-18",3,$null,$null,"2"),
    @("CO","KTLT",0,2,"0,1,2","8-12,12-18",2,$null,$null,"3,4"),
    @("CO","NMDT",0,2,"0,1,2","7-12,12-18",2,$null,$null,"3,4"),
    @("EE","MDT",0,2,"0,1,2","7-12,12-18",2,$null,$null,"5,6"),
    @("EE","VT",0,2,"0,1,2","7-12,12-18",2,$null,$null,"5"),
    @("CH","HDC",0,2,"0,1,2","7-12,12-18",2,$null,$null,"7"),
    @("CH","HHC",0,2,"0,1,2","7-12,12-18",2,$null,$null,"8")
)

for ($r = 0; $r -lt $subjectRows.Length; $r++) {
    $rowVals = $subjectRows[$r]
    for ($col = 0; $col -lt $rowVals.Length; $col++) {
        $val = $rowVals[$col]
        if ($null -ne $val) {
            $subject.Cells.Item($r + 2, $col + 1).Value = $val
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Insert a SUBJECT_ID column into "Enrollment" (becomes column C, pushing
#    the existing STUDENT_ID column to D) and fill it in.
# ---------------------------------------------------------------------------
$enrollment = $wb.Worksheets.Item("Enrollment")
$enrollment.Columns.Item(3).Insert() | Out-Null
$enrollment.Cells.Item(1, 3).Value = "SUBJECT_ID"

$subjectIds = @(0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,2,2,2,2,2,2,3,3,3,3,3,3,4,4,4,4,4,4,5,5,5,5,5,5,6,6,6,6,6,6,7,7,7,7,7,7)
for ($i = 0; $i -lt $subjectIds.Length; $i++) {
    $enrollment.Cells.Item($i + 2, 3).Value = $subjectIds[$i]
}

# ---------------------------------------------------------------------------
# 3. Drop the "NewSubject" sheet now that its data lives on "Subject".
# ---------------------------------------------------------------------------
$newSubject = $wb.Worksheets.Item("NewSubject")
$newSubject.Delete() | Out-Null
